{"js": "// Update the benchmark statistics table to reflect the corrected values.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of rowIndex -> new cell text (column 0, the only column in the table).\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"3600\",\n  4: \"0.00001\",\n  5: \"0.00262\",\n  6: \"0.00018\",\n  7: \"0.00005\",\n  8: \"0.00032\",\n  9: \"0.00039\",\n  10: \"0.00046\",\n  11: \"0.78965\",\n  43: \"99.85\",\n  44: \"0.79\",\n  45: \"514\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(parseInt(rowIndex, 10), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark statistics table to reflect the corrected values.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Map of 1-based row index -> new cell text (column 1, the only column).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"3600\"\n    5  = \"0.00001\"\n    6  = \"0.00262\"\n    7  = \"0.00018\"\n    8  = \"0.00005\"\n    9  = \"0.00032\"\n    10 = \"0.00039\"\n    11 = \"0.00046\"\n    12 = \"0.78965\"\n    44 = \"99.85\"\n    45 = \"0.79\"\n    46 = \"514\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
